$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 11871.429
$ws.Range("I69").Value = 11350
$ws.Range("J69").Value = 15000
$ws.Range("K69").Value = 34050
$ws.Range("L69").Value = 45000
$ws.Range("M69").Value = -33176
$ws.Range("N69").Value = -46748

$ws.Range("H72").Value = 11871.429
$ws.Range("I72").Value = 11350
$ws.Range("J72").Value = 15000
$ws.Range("K72").Value = 102150
$ws.Range("L72").Value = 135000
$ws.Range("M72").Value = -97782
$ws.Range("N72").Value = -143736.003

$ws.Range("H96").Value = 621.0769
$ws.Range("I96").Value = 573
$ws.Range("J96").Value = 662.2857
$ws.Range("K96").Value = 1719
$ws.Range("L96").Value = 1986.8571
$ws.Range("M96").Value = -346
$ws.Range("N96").Value = -4732.8571

$ws.Range("H98").Value = 2786.6843
$ws.Range("J98").Value = 2957.2856
$ws.Range("L98").Value = 2957.2856
$ws.Range("N98").Value = -5953.2856

$ws.Range("H100").Value = 1699
$ws.Range("I100").Value = 1344.2222
$ws.Range("J100").Value = 2497.25
$ws.Range("K100").Value = 1344.2222
$ws.Range("L100").Value = 2497.25
$ws.Range("M100").Value = -803.2221999999999
$ws.Range("N100").Value = -3579.25

$ws.Range("H101").Value = 1970.579
$ws.Range("I101").Value = 407
$ws.Range("K101").Value = 1221
$ws.Range("M101").Value = 401

$ws.Range("H122").Value = 2786.6843
$ws.Range("J122").Value = 2957.2856
$ws.Range("L122").Value = 8871.856800000001
$ws.Range("N122").Value = -13771.8568

$ws.Range("H132").Value = 3509.5
$ws.Range("I132").Value = 3250.68
$ws.Range("K132").Value = 9752.039999999999
$ws.Range("M132").Value = -7222.039999999999

$ws.Range("H138").Value = 4400.886
$ws.Range("J138").Value = 4608.183
$ws.Range("L138").Value = 13824.549
$ws.Range("N138").Value = -24104.549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1136.9014
$ws.Range("I32").Value = 738.0154
$ws.Range("K32").Value = 738.0154
$ws.Range("M32").Value = -451.0154

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H105").Value = 64500
$ws.Range("J105").Value = 64500
$ws.Range("L105").Value = 64500
$ws.Range("N105").Value = -71488

$ws.Range("H110").Value = 853.9
$ws.Range("I110").Value = 796.3570999999999
$ws.Range("K110").Value = 796.3570999999999
$ws.Range("M110").Value = 1248.6429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 362.77777
$ws.Range("J64").Value = 338.41666
$ws.Range("L64").Value = 338.41666
$ws.Range("N64").Value = -788.41666

$ws.Range("H67").Value = 362.77777
$ws.Range("J67").Value = 338.41666
$ws.Range("L67").Value = 338.41666
$ws.Range("N67").Value = -1898.41666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 13832.833
$ws.Range("J41").Value = 22500
$ws.Range("L41").Value = 22500
$ws.Range("N41").Value = -23356

$ws.Range("H51").Value = 54999
$ws.Range("J51").Value = 54999
$ws.Range("L51").Value = 54999
$ws.Range("N51").Value = -56471

$ws.Range("H60").Value = 21441.715
$ws.Range("J60").Value = 26999.8
$ws.Range("L60").Value = 26999.8
$ws.Range("N60").Value = -28021.8

$ws.Range("H61").Value = 54999
$ws.Range("J61").Value = 54999
$ws.Range("L61").Value = 54999
$ws.Range("N61").Value = -55695

$ws.Range("H137").Value = 34666.332
$ws.Range("J137").Value = 34666.332
$ws.Range("L137").Value = 34666.332
$ws.Range("N137").Value = -44866.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 955.9666999999999
$ws.Range("I113").Value = 863.7857
$ws.Range("J113").Value = 1036.625
$ws.Range("K113").Value = 2591.3571
$ws.Range("L113").Value = 3109.875
$ws.Range("M113").Value = -421.3571000000002
$ws.Range("N113").Value = -7449.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6988.2354

$ws.Range("H126").Value = 5172.591
$ws.Range("I126").Value = 3800.5334
$ws.Range("K126").Value = 11401.6002
$ws.Range("M126").Value = -8931.600199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 12508333
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 12508333
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 12508333
$ws.Range("N20").Value = -12508785
$ws.Range("M20").ClearContents()

$ws.Range("H24").Value = 501750
$ws.Range("I24").Value = 3500
$ws.Range("J24").Value = 1000000
$ws.Range("K24").Value = 3500
$ws.Range("L24").Value = 1000000
$ws.Range("M24").Value = -3157
$ws.Range("N24").Value = -1000686

$ws.Range("H68").Value = 2522.125
$ws.Range("J68").Value = 2442
$ws.Range("L68").Value = 2442
$ws.Range("N68").Value = -3940

$ws.Range("H71").Value = 2522.125
$ws.Range("J71").Value = 2442
$ws.Range("L71").Value = 12210
$ws.Range("N71").Value = -19698

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 6999.5
$ws.Range("I7").Value = 6999.5
$ws.Range("K7").Value = 6999.5
$ws.Range("M7").Value = -6886.5

$ws.Range("H37").Value = 20000
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("N37").Value = -20406

$ws.Range("H46").Value = 150000
$ws.Range("J46").Value = 150000
$ws.Range("L46").Value = 150000
$ws.Range("N46").Value = -150462

$ws.Range("H100").Value = 3133.625
$ws.Range("I100").Value = 315.25
$ws.Range("J100").Value = 5952
$ws.Range("K100").Value = 630.5
$ws.Range("L100").Value = 11904
$ws.Range("M100").Value = -89.5
$ws.Range("N100").Value = -12986

$ws.Range("H104").Value = 21149.334
$ws.Range("J104").Value = 21149.334
$ws.Range("L104").Value = 21149.334
$ws.Range("N104").Value = -28137.334

$ws.Range("H126").Value = 1468.4445
$ws.Range("J126").Value = 1004.3333
$ws.Range("L126").Value = 3012.9999
$ws.Range("N126").Value = -7952.9999

$ws.Range("H134").Value = 150000
$ws.Range("J134").Value = 150000
$ws.Range("L134").Value = 450000
$ws.Range("N134").Value = -455070
